# The deck ships two theme parts:
#   theme1.xml -> bound to the (only) slide master, currently the "Integral" theme
#   theme2.xml -> bound to the notes master, currently a generic "Office Theme"
#
# The commit swaps the two themes' content: the slide master's theme becomes
# the generic "Office Theme" colors, and the notes master's theme becomes the
# "Integral" colors. Font scheme / format scheme are identical between the two
# themes already, so only the 12 color-scheme entries (and the theme names)
# actually change.
#
# PowerPoint's automation model exposes the modern 12-slot theme color scheme
# through Slide.ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink, in that order). We use it against the presentation's (slide
# master) theme to push in the "Office Theme" colors that replace the
# current "Integral" ones.

$p = $ppt.ActivePresentation

# --- Helper: RGB() style packer (VBA/PowerPoint uses 0xBBGGRR ordering) ----
function ThemeRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Slide master's theme (theme1.xml): Integral colors -> Office colors
$slide = $p.Slides.Item(1)
$masterColors = $slide.ThemeColorScheme

$officeThemeColors = @(
    (ThemeRGB 0x00 0x00 0x00),   # dk1
    (ThemeRGB 0xFF 0xFF 0xFF),   # lt1
    (ThemeRGB 0x44 0x54 0x6A),   # dk2
    (ThemeRGB 0xE7 0xE6 0xE6),   # lt2
    (ThemeRGB 0x5B 0x9B 0xD5),   # accent1
    (ThemeRGB 0xED 0x7D 0x31),   # accent2
    (ThemeRGB 0xA5 0xA5 0xA5),   # accent3
    (ThemeRGB 0xFF 0xC0 0x00),   # accent4
    (ThemeRGB 0x44 0x72 0xC4),   # accent5
    (ThemeRGB 0x70 0xAD 0x47),   # accent6
    (ThemeRGB 0x05 0x63 0xC1),   # hlink
    (ThemeRGB 0x95 0x4F 0x72)    # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $masterColors.Colors($i + 1).RGB = $officeThemeColors[$i]
}

Write-Host "Slide master theme colors updated (Integral -> Office Theme)."
